$d = $word.ActiveDocument

# -----------------------------------------------------------------
# 1) Remove the stray "_GoBack" bookmark left inside the
#    "TABELAS (8)" heading paragraph.
# -----------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# -----------------------------------------------------------------
# 2) Italicise the "TAREFAFUNCIONARIO" bullet item (both the run
#    text and the paragraph mark get <w:i/><w:iCs/>).
# -----------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "TAREFAFUNCIONARIO`r") {
        $xml = $p.Range.WordOpenXML
        $plain = '<w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr>'
        $italic = '<w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:i/><w:iCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr>'
        $xml = $xml.Replace($plain, $italic)
        $p.Range.InsertXML($xml)
    }
}

# -----------------------------------------------------------------
# 3) Append the new "TUTORIAL" section (page break, heading, the
#    instructions paragraph and the three numbered / highlighted
#    steps, plus two trailing blank paragraphs) right after the
#    existing "CARGO" bullet, before the section properties.
# -----------------------------------------------------------------
$newParagraphsXml = @'
<w:p><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:br w:type="page"/></w:r></w:p><w:p><w:pPr><w:tabs><w:tab w:val="left" w:pos="6840"/></w:tabs><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:lastRenderedPageBreak/><w:t>TUTORIAL</w:t></w:r></w:p><w:p><w:pPr><w:tabs><w:tab w:val="left" w:pos="6840"/></w:tabs><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:tabs><w:tab w:val="left" w:pos="6840"/></w:tabs><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">1) Verificar se existe as Basedado, </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>IPGFuncionariosDbContext</w:t></w:r><w:r><w:t xml:space="preserve"> ou </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>ApplicationDbContext</w:t></w:r><w:r><w:t>, se no caso existir, apague.</w:t></w:r></w:p><w:p><w:pPr><w:tabs><w:tab w:val="left" w:pos="6840"/></w:tabs><w:rPr><w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:color w:val="FFFFFF"/><w:sz w:val="23"/><w:szCs w:val="23"/><w:shd w:val="clear" w:color="auto" w:fill="0099FF"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:color w:val="FFFFFF"/><w:sz w:val="23"/><w:szCs w:val="23"/><w:shd w:val="clear" w:color="auto" w:fill="0099FF"/></w:rPr><w:t>2</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:color w:val="FFFFFF"/><w:sz w:val="23"/><w:szCs w:val="23"/><w:shd w:val="clear" w:color="auto" w:fill="0099FF"/></w:rPr><w:t xml:space="preserve">) Update-Database -Context IPGFuncionariosDbContext </w:t></w:r></w:p><w:p><w:pPr><w:tabs><w:tab w:val="left" w:pos="6840"/></w:tabs><w:rPr><w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:color w:val="FFFFFF"/><w:sz w:val="23"/><w:szCs w:val="23"/><w:shd w:val="clear" w:color="auto" w:fill="0099FF"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:color w:val="FFFFFF"/><w:sz w:val="23"/><w:szCs w:val="23"/><w:shd w:val="clear" w:color="auto" w:fill="0099FF"/></w:rPr><w:t>3</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:color w:val="FFFFFF"/><w:sz w:val="23"/><w:szCs w:val="23"/><w:shd w:val="clear" w:color="auto" w:fill="0099FF"/></w:rPr><w:t xml:space="preserve">) </w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:r><w:rPr><w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:color w:val="FFFFFF"/><w:sz w:val="23"/><w:szCs w:val="23"/><w:shd w:val="clear" w:color="auto" w:fill="0099FF"/></w:rPr><w:t>Update-Database -Context ApplicationDbContext</w:t></w:r><w:bookmarkEnd w:id="0"/></w:p><w:p><w:pPr><w:tabs><w:tab w:val="left" w:pos="6840"/></w:tabs><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:tabs><w:tab w:val="left" w:pos="6840"/></w:tabs><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr></w:p>
'@

$pkg = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">PKG_BODY_PLACEHOLDER</w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$pkg = $pkg.Replace("PKG_BODY_PLACEHOLDER", "<w:body>" + $newParagraphsXml + "</w:body>")

$endRange = $d.Range($d.Content.End, $d.Content.End)
$endRange.InsertXML($pkg)

Write-Output "done"
